# Apply corrections to the "listado de control" and "migracion de datos" rows
# in Hoja1, and mark "Reporte para contador" as "en proceso".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 63: "Error en listado de control - pagos de la fecha" -> mark done (100%)
$ws.Range("C63").Value = 1
$ws.Range("C63").NumberFormat = "0%"

# Row 64: "Migracion de datos" -> mark done (100%)
$ws.Range("C64").Value = 1
$ws.Range("C64").NumberFormat = "0%"

# Row 67: "Reporte para contador" -> mark as "en proceso"
$ws.Range("C67").Value = "en proceso"

# Update the view to reflect where the user ended up looking
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("C68").Select()
